# Update for release to deploy 0.1.1
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Metadata"
$ws2 = $wb.Worksheets.Item(2)   # "Include from NMDP Practitione"

# 1. Rename the "Include" sheet.
$ws2.Name = "Include #0"

# 2. Bump the Version and Date metadata values.
$ws1.Range("B3").Value = "0.1.1"
$ws1.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# 3. Insert a new "Jurisdiction" metadata row right after "Contact" (row 10).
#    This pushes Description/Purpose/Copyright/Immutable down one row
#    (now rows 12-15). Copy formatting from the row above so the new row
#    keeps the same border/wrap/style as the rest of the table.
$ws1.Rows.Item(11).Insert()
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""
